$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1202.8
$ws.Range("I2").Value = 138.33333
$ws.Range("K2").Value = 138.33333
$ws.Range("M2").Value = -25.33332999999999
$ws.Range("H11").Value = 790.5
$ws.Range("I11").Value = 790.5
$ws.Range("K11").Value = 790.5
$ws.Range("M11").Value = -650.5
$ws.Range("H19").Value = 280.2857
$ws.Range("J19").Value = 204
$ws.Range("L19").Value = 204
$ws.Range("N19").Value = -554
$ws.Range("H40").Value = 2493.6667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2493.6667
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2493.6667
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2843.6667
$ws.Range("H53").Value = 566.4545000000001
$ws.Range("I53").Value = 831.8333
$ws.Range("J53").Value = 248
$ws.Range("K53").Value = 831.8333
$ws.Range("L53").Value = 248
$ws.Range("M53").Value = -194.8333
$ws.Range("N53").Value = -1522
$ws.Range("H74").Value = 3659.6667
$ws.Range("I74").Value = 3989.5
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3989.5
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -3053.5
$ws.Range("N74").Value = -4872
$ws.Range("H77").Value = 3659.6667
$ws.Range("I77").Value = 3989.5
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 19947.5
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -15267.5
$ws.Range("N77").Value = -24360
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("H127").Value = 5640
$ws.Range("I127").Value = 3733.3333
$ws.Range("J127").Value = 8500
$ws.Range("K127").Value = 11199.9999
$ws.Range("L127").Value = 25500
$ws.Range("M127").Value = -6239.999899999999
$ws.Range("N127").Value = -35420
$ws.Range("H138").Value = 1732.4783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2099.6
$ws.Range("J61").Value = 2199.5
$ws.Range("L61").Value = 2199.5
$ws.Range("N61").Value = -2623.5
$ws.Range("H76").Value = 60000
$ws.Range("J76").Value = 60000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60676
$ws.Range("H79").Value = 60000
$ws.Range("J79").Value = 60000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62340
$ws.Range("H136").Value = 2099.6
$ws.Range("J136").Value = 2199.5
$ws.Range("L136").Value = 6598.5
$ws.Range("N136").Value = -11698.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 220
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 200
$ws.Range("M7").Value = -87
$ws.Range("H105").Value = 1956.625
$ws.Range("I105").Value = 1500.5
$ws.Range("K105").Value = 1500.5
$ws.Range("M105").Value = 246.5
$ws.Range("H134").Value = 10400
$ws.Range("I134").Value = 10457.143
$ws.Range("K134").Value = 31371.429
$ws.Range("M134").Value = -28836.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2683
$ws.Range("I134").Value = 1099.7
$ws.Range("K134").Value = 3299.1
$ws.Range("M134").Value = -764.1000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 408.5
$ws.Range("I2").Value = 199
$ws.Range("J2").Value = 431.77777
$ws.Range("K2").Value = 1194
$ws.Range("L2").Value = 2590.66662
$ws.Range("M2").Value = -1081
$ws.Range("N2").Value = -2816.66662
$ws.Range("H4").Value = 2352.647
$ws.Range("J4").Value = 2366.111
$ws.Range("L4").Value = 7098.333
$ws.Range("N4").Value = -7322.333
$ws.Range("H37").Value = 97498.336
$ws.Range("J37").Value = 97498.336
$ws.Range("L37").Value = 292495.008
$ws.Range("N37").Value = -292719.008
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5064
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13320
$ws.Range("H107").Value = 1603.3846
$ws.Range("I107").Value = 1199.8334
$ws.Range("J107").Value = 1949.2858
$ws.Range("K107").Value = 3599.5002
$ws.Range("L107").Value = 5847.857400000001
$ws.Range("M107").Value = -1679.5002
$ws.Range("N107").Value = -9687.857400000001
$ws.Range("H131").Value = 2245.25
$ws.Range("I131").Value = 1009
$ws.Range("J131").Value = 2492.5
$ws.Range("K131").Value = 3027
$ws.Range("L131").Value = 7477.5
$ws.Range("M131").Value = 2013
$ws.Range("N131").Value = -17557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 44638
$ws.Range("J20").Value = 44638
$ws.Range("L20").Value = 44638
$ws.Range("N20").Value = -45128
$ws.Range("H49").Value = 2075
$ws.Range("J49").Value = 2075
$ws.Range("L49").Value = 2075
$ws.Range("N49").Value = -2443
$ws.Range("H102").Value = 635.8823
$ws.Range("I102").Value = 613.125
$ws.Range("K102").Value = 613.125
$ws.Range("M102").Value = 1008.875
$ws.Range("H122").Value = 11371500
$ws.Range("I122").Value = 25011500
$ws.Range("K122").Value = 75034500
$ws.Range("M122").Value = -75032050
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1578.909
$ws.Range("I16").Value = 1640.8889
$ws.Range("K16").Value = 1640.8889
$ws.Range("M16").Value = -1470.8889
$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714
$ws.Range("H46").Value = 3315.2307
$ws.Range("J46").Value = 4285.4287
$ws.Range("L46").Value = 4285.4287
$ws.Range("N46").Value = -4661.4287
$ws.Range("H82").Value = 1945
$ws.Range("H85").Value = 1945
$ws.Range("H132").Value = 7233.148
$ws.Range("J132").Value = 8333
$ws.Range("L132").Value = 24999
$ws.Range("N132").Value = -30059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 357.7143
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 357.7143
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H126").Value = 3867.0715
$ws.Range("I126").Value = 3614
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 10842
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -8372
$ws.Range("N126").Value = -18439.25
$ws.Range("H132").Value = 1123.7273
$ws.Range("I132").Value = 1123.7273
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3371.1819
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -841.1819
$ws.Range("N132").ClearContents()
